$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (some Price values look
# like plain numbers, e.g. "0.9978" or "240.64", and would otherwise be
# auto-converted to numeric cells by Excel). Resetting the Style to "Normal"
# after the write keeps the original (default) cell formatting.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Price (D) / Volume(1h) (E) refresh for rows 2-46 ---
Set-TextValue $ws.Range("D2") '29.383.17'
Set-TextValue $ws.Range("E2") '  -0.03%  '
Set-TextValue $ws.Range("D3") '1.847.12'
Set-TextValue $ws.Range("E3") '  -0.10%  '
Set-TextValue $ws.Range("D4") '0.9978'
Set-TextValue $ws.Range("E4") '  -0.24%  '
Set-TextValue $ws.Range("D5") '240.64'
Set-TextValue $ws.Range("E5") '  +0.09%  '
Set-TextValue $ws.Range("D6") '0.6307'
Set-TextValue $ws.Range("E6") '  +0.51%  '
Set-TextValue $ws.Range("E7") '  -0.17%  '
Set-TextValue $ws.Range("D8") '0.07507'
Set-TextValue $ws.Range("E8") '  -1.67%  '
Set-TextValue $ws.Range("E9") '  -0.08%  '
Set-TextValue $ws.Range("D10") '24.44'
Set-TextValue $ws.Range("E10") '  -1.16%  '
Set-TextValue $ws.Range("D11") '0.07710'
Set-TextValue $ws.Range("E11") '  -0.45%  '
Set-TextValue $ws.Range("D12") '1.846.08'
Set-TextValue $ws.Range("E12") '  -2.23%  '
Set-TextValue $ws.Range("D13") '5.008'
Set-TextValue $ws.Range("E13") '  -0.49%  '
Set-TextValue $ws.Range("D14") '0.6810'
Set-TextValue $ws.Range("E14") '  +0.38%  '
Set-TextValue $ws.Range("D15") '0.00001031'
Set-TextValue $ws.Range("E15") '  -2.93%  '
Set-TextValue $ws.Range("D16") '82.21'
Set-TextValue $ws.Range("E16") '  -1.21%  '
Set-TextValue $ws.Range("D17") '2.103.68'
Set-TextValue $ws.Range("E17") '  -3.81%  '
Set-TextValue $ws.Range("D18") '6.160'
Set-TextValue $ws.Range("E18") '  +0.12%  '
Set-TextValue $ws.Range("D19") '29.387.79'
Set-TextValue $ws.Range("E19") '  -0.13%  '
Set-TextValue $ws.Range("D20") '229.70'
Set-TextValue $ws.Range("E20") '  +1.33%  '
Set-TextValue $ws.Range("D21") '12.35'
Set-TextValue $ws.Range("E21") '  +0.15%  '
Set-TextValue $ws.Range("D22") '0.9994'
Set-TextValue $ws.Range("E22") '  -0.14%  '
Set-TextValue $ws.Range("D23") '7.451'
Set-TextValue $ws.Range("E23") '  -0.57%  '
Set-TextValue $ws.Range("D24") '0.9991'
Set-TextValue $ws.Range("E24") '  -0.26%  '
Set-TextValue $ws.Range("D25") '158.88'
Set-TextValue $ws.Range("E25") '  +0.58%  '
Set-TextValue $ws.Range("D26") '0.1381'
Set-TextValue $ws.Range("E26") '  +0.13%  '
Set-TextValue $ws.Range("D27") '8.418'
Set-TextValue $ws.Range("E27") '  +0.01%  '
Set-TextValue $ws.Range("D28") '17.58'
Set-TextValue $ws.Range("E28") '  -0.55%  '
Set-TextValue $ws.Range("D29") '0.06413'
Set-TextValue $ws.Range("E29") '  +14.83%  '
Set-TextValue $ws.Range("D30") '1.387'
Set-TextValue $ws.Range("E30") '  +0.32%  '
Set-TextValue $ws.Range("D31") '1.474'
Set-TextValue $ws.Range("E31") '  +0.62%  '
Set-TextValue $ws.Range("D32") '4.094'
Set-TextValue $ws.Range("E32") '  -0.80%  '
Set-TextValue $ws.Range("D33") '4.060'
Set-TextValue $ws.Range("E33") '  +0.18%  '
Set-TextValue $ws.Range("D34") '1.819'
Set-TextValue $ws.Range("E34") '  -0.89%  '
Set-TextValue $ws.Range("D35") '1.143'
Set-TextValue $ws.Range("E35") '  -1.70%  '
Set-TextValue $ws.Range("D36") '0.6967'
Set-TextValue $ws.Range("E36") '  +0.10%  '
Set-TextValue $ws.Range("D37") '2.578'
Set-TextValue $ws.Range("E37") '  -0.41%  '
Set-TextValue $ws.Range("D38") '1.257.61'
Set-TextValue $ws.Range("E38") '  +2.22%  '
Set-TextValue $ws.Range("D39") '2.835'
Set-TextValue $ws.Range("E39") '  +4.34%  '
Set-TextValue $ws.Range("D40") '0.01825'
Set-TextValue $ws.Range("E40") '  +1.34%  '
Set-TextValue $ws.Range("E41") '  +2.84%  '
Set-TextValue $ws.Range("D42") '0.9079'
Set-TextValue $ws.Range("E42") '  +0.50%  '
Set-TextValue $ws.Range("D43") '0.9986'
Set-TextValue $ws.Range("E43") '  -0.29%  '
Set-TextValue $ws.Range("D44") '2.006.20'
Set-TextValue $ws.Range("E44") '  -18.47%  '
Set-TextValue $ws.Range("D45") '101.37'
Set-TextValue $ws.Range("E45") '  -0.23%  '
Set-TextValue $ws.Range("D46") '66.30'
Set-TextValue $ws.Range("E46") '  +0.60%  '

# --- Rows 47-51: BabyDogeCoin dropped out of the table, every row below it
# shifted up by one, and TheSandbox was newly added as row 51 ---
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D47") '0.1177'
Set-TextValue $ws.Range("E47") '  +2.87%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D48") '7.065'
Set-TextValue $ws.Range("E48") '  -1.48%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D49") '1.712'
Set-TextValue $ws.Range("E49") '  +2.16%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D50") '9.023'
Set-TextValue $ws.Range("E50") '  +0.59%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D51") '0.3944'
Set-TextValue $ws.Range("E51") '  -1.65%  '
